# Update the two-digit multiplication problems in the document.
# Each "AxB=" string in the document is unique, so a straightforward
# Find/Replace (no wildcards) for each pair is safe and unambiguous.

$d = $word.ActiveDocument

$replacements = @(
    @("91×42=", "54×51="),
    @("15×12=", "54×37="),
    @("19×54=", "33×85="),
    @("26×18=", "70×65="),
    @("58×76=", "35×20="),
    @("57×20=", "23×18="),
    @("21×82=", "52×13="),
    @("33×61=", "52×99="),
    @("45×99=", "71×28="),
    @("64×92=", "17×55="),
    @("65×23=", "22×74="),
    @("41×21=", "65×70="),
    @("41×39=", "91×12="),
    @("37×46=", "66×76="),
    @("93×69=", "16×18="),
    @("48×73=", "65×36="),
    @("83×37=", "19×90="),
    @("24×56=", "48×31="),
    @("27×63=", "16×94="),
    @("22×51=", "60×12="),
    @("70×59=", "47×21="),
    @("60×97=", "78×89="),
    @("54×98=", "44×17="),
    @("11×99=", "79×59="),
    @("45×49=", "85×91=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
